{"js": "// Apply targeted text replacements for the Butterfly Pavilion caption update.\nconst replacements = [\n  [\"Explore a magical butterfly garden indoors with hundreds of fluttering butterflies, lush native plants, and natural light. See all life stages\u2014from eggs to chrysalises\u2014and chat with educators in this enchanting experience perfect for families.\", \"Step into a magical garden at NHM\u2019s Butterfly Pavilion and witness daily butterfly flights -- from caterpillars to chrysalis to airborne beauties. Perfectly timed 30-minute visits offer up-close nature moments for all ages.\"],\n  [\"\ud83d\udccd Location: Natural History Museum of LA County\", \"\ud83d\udccd Location: Natural History Museum\"],\n  [\"\ud83d\udcc5 Date: 2025\u201107\u201110\", \"\ud83d\udcc5 Date: 2025-07-17 \u2013 2025-07-20\"],\n  [\"\ud83d\udd58 Time: 9:30 AM\u20135 PM\", \"\ud83d\udd58 Time: 10:00 am \u2013 4:30 pm reservation slots\"],\n  [\"\ud83d\udcb0 Tickets: $10 + general admission\", \"\ud83d\udcb0 Tickets: $10 + general admission; Members free\"],\n  [\"#ButterflyMagic #NHMLA #FamilySTEM #NatureExhibit #InsectLife #LAKids #InteractiveLearning #SeasonalExhibit #SummerOuting #ShitToDoWithKids #shittodowithkids #stdwkids #familyactivities #kidslosangeles\", \"#ButterflyPavilion #NHMLA #LiveButterflies #FamilyOuting #InteractiveNature #ScienceFun #LAEvents #SeasonalExhibit #ReservationRequired #ShitToDoWithKids #shittodowithkids #stdwkids #familyactivities #kidslosangeles\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n  if (found.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  found.items[0].insertText(newText, \"Replace\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n}\n", "ps1": "# Update the Butterfly Pavilion social caption body text.\n# The whole caption lives in a single paragraph/run, with the various\n# lines separated by manual line breaks (vertical-tab / <w:br/>), so each\n# logical \"line\" is replaced in place via Find/Replace on $d.Content.\n$replacements = @(\n    @{ Old = 'Explore a magical butterfly garden indoors with hundreds of fluttering butterflies, lush native plants, and natural light. See all life stages\u2014from eggs to chrysalises\u2014and chat with educators in this enchanting experience perfect for families.'; New = 'Step into a magical garden at NHM\u2019s Butterfly Pavilion and witness daily butterfly flights -- from caterpillars to chrysalis to airborne beauties. Perfectly timed 30-minute visits offer up-close nature moments for all ages.' }\n    @{ Old = '\ud83d\udccd Location: Natural History Museum of LA County'; New = '\ud83d\udccd Location: Natural History Museum' }\n    @{ Old = '\ud83d\udcc5 Date: 2025\u201107\u201110'; New = '\ud83d\udcc5 Date: 2025-07-17 \u2013 2025-07-20' }\n    @{ Old = '\ud83d\udd58 Time: 9:30 AM\u20135 PM'; New = '\ud83d\udd58 Time: 10:00 am \u2013 4:30 pm reservation slots' }\n    @{ Old = '\ud83d\udcb0 Tickets: $10 + general admission'; New = '\ud83d\udcb0 Tickets: $10 + general admission; Members free' }\n    @{ Old = '#ButterflyMagic #NHMLA #FamilySTEM #NatureExhibit #InsectLife #LAKids #InteractiveLearning #SeasonalExhibit #SummerOuting #ShitToDoWithKids #shittodowithkids #stdwkids #familyactivities #kidslosangeles'; New = '#ButterflyPavilion #NHMLA #LiveButterflies #FamilyOuting #InteractiveNature #ScienceFun #LAEvents #SeasonalExhibit #ReservationRequired #ShitToDoWithKids #shittodowithkids #stdwkids #familyactivities #kidslosangeles' }\n)\n\n$d = $word.ActiveDocument\n\nforeach ($item in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $item.Old,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $item.New,   # ReplaceWith\n        1            # Replace (wdReplaceOne)\n    )\n    if (-not $found) {\n        throw \"edit.ps1: could not find expected text: $($item.Old)\"\n    }\n}\n"}
